# Apply the "1st Commit on 3rd May" edit:
#  - TestCases!D2 = "PASS", D3 = "PASS", D4 = "FAIL"
#  - TestSteps!H2:H33 = "PASS" for every data row
#  - TestSteps becomes the active/selected sheet (was TestCases)

$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("TestCases")
$wsSteps = $wb.Worksheets.Item("TestSteps")

# --- TestCases sheet: add Results column values ---
$wsCases.Range("D2").Value = "PASS"
$wsCases.Range("D3").Value = "PASS"
$wsCases.Range("D4").Value = "FAIL"

# --- TestSteps sheet: add Results column values for every row ---
for ($r = 2; $r -le 33; $r++) {
    $wsSteps.Cells.Item($r, 8).Value = "PASS"
}

# --- Selection / active sheet bookkeeping ---
$null = $wsCases.Range("D2:D4").Select()
$null = $wsSteps.Range("G35").Select()
$wsSteps.Activate()
